$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append new trade row (row 13) below the existing data (rows 2-12)
$ws.Cells.Item(13, 1).Value = 9758.9
$ws.Cells.Item(13, 2).Value = 9865.4500000000007
$ws.Cells.Item(13, 3).Value = 281.06
$ws.Cells.Item(13, 4).Value = 284.08999999999997
$ws.Cells.Item(13, 5).Value = $true
$ws.Cells.Item(13, 6).Value = 1.08

# Copy the date/time format used for the rest of column G (style index 1, numFmtId 22)
# before assigning the value, so it reuses the existing style instead of creating a new one.
$ws.Cells.Item(12, 7).Copy()
$ws.Cells.Item(13, 7).PasteSpecial(-4122)
$ws.Cells.Item(13, 7).Value = 42620.766412037039

$ws.Cells.Item(13, 8).Value = $false
